$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column D ("is_normal_for_donor") to the right of the "n"
# column, shifting the old D ("test") and E ("test_sample.1.bam") columns
# to E and F respectively.
$ws.Columns("D:D").Insert()

# Match the inherited width of the neighbouring column (Excel carries the
# left-hand column's formatting onto a freshly inserted column).
$ws.Columns("D:D").ColumnWidth = $ws.Columns("C:C").ColumnWidth

# Header for the newly inserted column.
$ws.Range("D1").Value = "Y"

# Selection left where the editing session ended.
$ws.Range("D5").Select()
